$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark "Scramble MMSIs*" (row 10) as complete ---
$ws.Range("A10").Value = "x"

# --- Insert 3 new rows at position 13 (pushes old rows 13-16 down to 16-19) ---
$ws.Rows("13:15").Insert()

# Inherit the formatting pattern (styles 4,3,3,3 for A:D) from row 12 for the three new rows
$ws.Range("A12:D12").Copy()
$ws.Range("A13:D15").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Row 13: "Add spoof remover into vectorization script" ---
$ws.Range("C13").Value = "Add spoof remover into vectorization script"
$ws.Range("D13").Value = "Ended up not needing it because I implemented a 60 km limit between successive points. "
$ws.Range("A13").Value = "Not implemented"
$ws.Range("B13").Value = "Short"

# E13 note, formatted like the existing bordered note column (same style as D18/D19)
$ws.Range("D18:D18").Copy()
$ws.Range("E13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("E13").Value = "NEED TO TEST WHETHER THIS LIMIT BIASES THE DATA IN EARLY YEARS WITH LESS SATELLITE COVERAGE… "

$ws.Rows("13").RowHeight = 30

# --- Row 14: "Improved ship type assignment" ---
$ws.Range("C14").Value = "Improved ship type assignment"
$ws.Range("B14").Value = "Short"
$ws.Range("D14").Value = "Instead of just taking the metadata from the first static message transmitted for each ship each day, the new script takes the most frequently transmitted ship type in all static messages for each ship in each month (downweighting NA and 0 entries). This minimizes the overall chance of error. Also, the new script assigns static information from vessels to transit segments based on monthly transmissions rather than daily (AISlookup is joined to AISsf based on MMSI, not AIS_ID). "
$ws.Range("A14").Value = "x"

# E14 stays blank but picks up a wrap-text style
$ws.Range("E14").WrapText = $true

$ws.Rows("14").RowHeight = 75

# --- The "Short" label cell below the new rows (B18) now shares the same
#     bordered/shaded style as the note cell beside it (D18) ---
$ws.Range("D18:D18").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Column A is a bit wider now to fit "Not implemented" ---
$ws.Columns("A").ColumnWidth = 16.67

# --- Active cell selection as left by the author ---
$ws.Range("A15").Select()

Write-Output "done"
